# Nieuwe data toegevoegd via Streamlit op 2024-12-04 11:21:08
# Append a new data row (row 96) to the worksheet, mirroring the
# structure of the existing rows:
#   A: Bedrijfnaam, B: Vestiging, C: Type Opvangvoorziening,
#   D: Rapportdatum (stored as plain text, e.g. "2024-04-29"),
#   E-J: numeric score columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 96

# Column D holds a date-like string that must stay as literal text
# (not get auto-converted into a date serial number by Excel).
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "2024-04-29"
# Drop the temporary text-format style so the cell ends up unstyled,
# just like the other data cells in the sheet.
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 1).Value = "Kindergarden"
$ws.Cells.Item($row, 2).Value = "Kindergarden Den Haag Eisenhowerlaan"
$ws.Cells.Item($row, 3).Value = "KDV"
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
